$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet per the diff (Gamma2F-HW50.xpc -> Gamma2F)
$ws.Name = "Gamma2F"

# Append a new row 16, mirroring the pattern of row 15 (A=index, B=label, C:M=1)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:M16").Value = 1
